$wb = $excel.ActiveWorkbook

# --- REGCA1: remove Iqmax/Iqmin and Iqrmax/Iqrmin columns ---
$ws1 = $wb.Worksheets.Item("REGCA1")
$ws1.Range("U1:V2").Delete(-4159) # xlShiftToLeft, removes Iqmax/Iqmin
$ws1.Range("R1:S2").Delete(-4159) # xlShiftToLeft, removes Iqrmax/Iqrmin

# --- REPCA1: insert PLflag column after Fflag ---
$ws2 = $wb.Worksheets.Item("REPCA1")
$ws2.Range("L1:L2").Insert(-4161) # xlShiftToRight
$ws2.Range("L1").Value = "PLflag"
$ws2.Range("L2").Value = 0

Write-Host ("REGCA1 UsedRange=" + $ws1.UsedRange.Address())
Write-Host ("REPCA1 UsedRange=" + $ws2.UsedRange.Address())
